# "added update and delete project"
#
# Project 1 (row 2) is edited in place: the project name is changed from
# "aqsw" to "wwaa" and the (now stale) student 1 / student 2 details are
# cleared out ("delete").
#
# Project 2 (row 3) is a brand-new project ("added"): S no. 2, project
# name "aa", with its student detail columns left blank.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update project 1 -------------------------------------------------
$ws.Range("B2").Value = "wwaa"
$ws.Range("C2:H2").ClearContents()

# --- Add project 2 -----------------------------------------------------
$ws.Range("A3").Value = 2
$ws.Range("B3").Value = "aa"
